$wb = $excel.ActiveWorkbook

# Handback report generation: update the "Correspond Handoff Datetime" (E)
# and "Correspond Handback DateTime" (H) columns for the second handback
# file row (row 3) on each locale sheet.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-12 12:36:10"
$wsZhCn.Range("H3").Value = "2016-03-12 12:36:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-12 12:36:13"
$wsDeDe.Range("H3").Value = "2016-03-12 12:36:32"
